{"js": "// \"Remove VAR, add shorthand notation for INS/EXEC\"\n//\n// 1. The stray \"_GoBack\" bookmark that sits between the \"ALIAS\" run and the\n//    \" delta2 INS $wp.endMilestone.plannedDelta+++\" run is removed.\n// 2. That same bookmark re-appears between \"+++\" and \"FOR\" in the\n//    \"+++ FOR wp IN project.workPackages +++\" cell (note: the space that used\n//    to separate \"+++\" and \"FOR\" is not retained once the text is split).\n// 3. Every \"+++[NAME]+++\" placeholder (the \"VAR\" form) becomes the shorthand\n//    \"+++*NAME+++\" form (used for acronym/title/leader/starts/delta1/ends/delta2).\n\nconst body = context.document.body;\n\n// --- 1. Drop the old \"_GoBack\" bookmark wherever it currently lives. ---\ncontext.document.deleteBookmark(\"_GoBack\");\nawait context.sync();\n\n// --- 2. Turn \"+++ FOR\" into \"+++\" + bookmark(\"_GoBack\") + \"FOR\". ---\nconst forSearch = body.search(\"+++ FOR\", { matchCase: true });\nforSearch.load(\"items\");\nawait context.sync();\n\nif (forSearch.items.length > 0) {\n  const forRange = forSearch.items[0];\n  forRange.insertText(\"+++FOR\", Word.InsertLocation.replace);\n  await context.sync();\n\n  // Re-locate the freshly written text, then find the \"+++\" prefix inside it\n  // so the bookmark can be anchored exactly between \"+++\" and \"FOR\".\n  const reSearch = body.search(\"+++FOR\", { matchCase: true });\n  reSearch.load(\"items\");\n  await context.sync();\n\n  const wholeRange = reSearch.items[0];\n  const prefixSearch = wholeRange.search(\"+++\", { matchCase: true });\n  prefixSearch.load(\"items\");\n  await context.sync();\n\n  const prefixRange = prefixSearch.items[0];\n  const boundary = prefixRange.getRange(Word.RangeLocation.end);\n  boundary.insertBookmark(\"_GoBack\");\n  await context.sync();\n}\n\n// --- 3. Shorthand notation: \"+++[\" -> \"+++*\" and \"]+++\" -> \"+++\". ---\nconst openBrackets = body.search(\"+++[\", { matchCase: true });\nopenBrackets.load(\"items\");\nawait context.sync();\nfor (const r of openBrackets.items) {\n  r.insertText(\"+++*\", Word.InsertLocation.replace);\n}\nawait context.sync();\n\nconst closeBrackets = body.search(\"]+++\", { matchCase: true });\ncloseBrackets.load(\"items\");\nawait context.sync();\nfor (const r of closeBrackets.items) {\n  r.insertText(\"+++\", Word.InsertLocation.replace);\n}\nawait context.sync();\n", "ps1": "# \"Remove VAR, add shorthand notation for INS/EXEC\"\n#\n# 1. The stray \"_GoBack\" bookmark that currently sits between the \"ALIAS\" run\n#    and the \" delta2 INS $wp.endMilestone.plannedDelta+++\" run is removed.\n# 2. That same bookmark re-appears between \"+++\" and \"FOR\" in the\n#    \"+++ FOR wp IN project.workPackages +++\" cell (the space that used to\n#    separate \"+++\" and \"FOR\" is not retained once the text is split).\n# 3. Every \"+++[NAME]+++\" placeholder (the \"VAR\" form) becomes the shorthand\n#    \"+++*NAME+++\" form (acronym/title/leader/starts/delta1/ends/delta2).\n\n$d = $word.ActiveDocument\n\n# --- 1. Drop the old \"_GoBack\" bookmark wherever it currently lives. ---\nif ($d.Bookmarks.Exists(\"_GoBack\")) {\n    $d.Bookmarks.Item(\"_GoBack\").Delete()\n}\n\n# --- 2. Turn \"+++ FOR\" into \"+++\" + bookmark(\"_GoBack\") + \"FOR\". ---\n$forRange = $d.Content\nif ($forRange.Find.Execute(\"+++ FOR\")) {\n    $forRange.Text = \"+++FOR\"\n    $boundary = $forRange.Start + 3\n    $bmRange = $d.Range($boundary, $boundary)\n    $d.Bookmarks.Add(\"_GoBack\", $bmRange)\n}\n\n# --- 3. Shorthand notation: \"+++[\" -> \"+++*\" and \"]+++\" -> \"+++\". ---\n$openRange = $d.Content\n$openRange.Find.Execute(\"+++[\", $false, $false, $false, $false, $false, $true, 1, $false, \"+++*\", 2)\n\n$closeRange = $d.Content\n$closeRange.Find.Execute(\"]+++\", $false, $false, $false, $false, $false, $true, 1, $false, \"+++\", 2)\n"}
